$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.256.56'
$ws.Range("E2").Value = '  +4.30%  '
$ws.Range("D3").Value = '2.493.37'
$ws.Range("E3").Value = '  +2.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.75'
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.08'
$ws.Range("E6").Value = '  +5.03%  '
$ws.Range("E7").Value = '  +2.24%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.542'
$ws.Range("E9").Value = '  +2.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.29'
$ws.Range("E10").Value = '  +7.87%  '
$ws.Range("E11").Value = '  +1.95%  '
$ws.Range("E12").Value = '  +1.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.50'
$ws.Range("E13").Value = '  +1.19%  '
$ws.Range("E14").Value = '  +2.31%  '
$ws.Range("D15").Value = '2.882.82'
$ws.Range("E15").Value = '  +2.76%  '
$ws.Range("D16").Value = '2.503.20'
$ws.Range("E16").Value = '  +2.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.857'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '47.191.33'
$ws.Range("E18").Value = '  +4.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.98'
$ws.Range("E19").Value = '  +6.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.69'
$ws.Range("E20").Value = '  +5.54%  '
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.76'
$ws.Range("E22").Value = '  +2.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.46'
$ws.Range("E23").Value = '  +7.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '250.58'
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("E25").Value = '  +4.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.31'
$ws.Range("E26").Value = '  +2.16%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.29'
$ws.Range("E28").Value = '  +1.56%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.03'
$ws.Range("E29").Value = '  +4.21%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.02'
$ws.Range("E30").Value = '  +6.43%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.138'
$ws.Range("E31").Value = '  +10.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.37'
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.49'
$ws.Range("E33").Value = '  +5.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.69'
$ws.Range("E34").Value = '  -2.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0793'
$ws.Range("E35").Value = '  +4.14%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.98'
$ws.Range("E37").Value = '  +5.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.70'
$ws.Range("E38").Value = '  +5.99%  '
$ws.Range("E39").Value = '  +3.97%  '
$ws.Range("E40").Value = '  +2.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '122.02'
$ws.Range("E41").Value = '  -4.16%  '
$ws.Range("E42").Value = '  +2.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.09'
$ws.Range("E43").Value = '  +2.31%  '
$ws.Range("E44").Value = '  +3.07%  '
$ws.Range("D45").Value = '1.967.33'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.01'
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.04'
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("E50").Value = '  +9.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.46'
$ws.Range("E51").Value = '  +3.82%  '
